$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card8")

$ws.Range("O1").Value = "Serviced by"

$ws.Range("O2").Value = "nan"
$ws.Range("O3").Value = "nan"
$ws.Range("O4").Value = "nan"
$ws.Range("O5").Value = "nan"
$ws.Range("O6").Value = "nan"
$ws.Range("O7").Value = "nan"
$ws.Range("O8").Value = "م.محمد عبدالله ،ف.محمود ايهاب"
$ws.Range("O9").Value = "nan"
$ws.Range("O10").Value = "nan"
$ws.Range("O11").Value = "nan"
$ws.Range("O12").Value = "nan"
$ws.Range("O13").Value = "nan"
